$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the ticker grid (B2:F22) with the latest scan results and append
# three new rows (23:25) continuing the numbered list in column A.

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "NSE:ARVSMART"
$ws.Range("C2").Value = "NSE:AGARIND"
$ws.Range("D2").Value = "NSE:LTIM"
$ws.Range("E2").Value = "NSE:BEL"
$ws.Range("F2").Value = "NSE:360ONE"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "NSE:BCG"
$ws.Range("C3").Value = "NSE:ETHOSLTD"
$ws.Range("D3").Value = "NSE:MARICO"
$ws.Range("E3").Value = "NSE:CIPLA"
$ws.Range("F3").Value = "NSE:CYIENT"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "NSE:CMSINFO"
$ws.Range("C4").Value = "NSE:HDFCGOLD"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "NSE:COALINDIA"
$ws.Range("F4").Value = ""

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "NSE:DWARKESH"
$ws.Range("C5").Value = "NSE:MOHITIND"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "NSE:HDFCBANK"
$ws.Range("F5").Value = ""

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "NSE:FACT"
$ws.Range("C6").Value = "NSE:NELCO"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "NSE:FILATEX"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "NSE:GNFC"
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "NSE:GSFC"
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "NSE:HESTERBIO"
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "NSE:HONDAPOWER"
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "NSE:ICRA"
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "NSE:IKIO"
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NSE:INFIBEAM"
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "NSE:KAUSHALYA"
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "NSE:LTTS"
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "NSE:MADRASFERT"
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "NSE:MAHEPC"
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "NSE:MOREPENLAB"
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "NSE:OSWALAGRO"
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "NSE:PALASHSECU"
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = ""

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "NSE:PENINLAND"
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = ""

# New rows appended at the bottom of the list
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "NSE:RAILTEL"
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = ""

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "NSE:RCF"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""
$ws.Range("E24").Value = ""
$ws.Range("F24").Value = ""

$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "NSE:ROHLTD"
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""

# Carry over the bold / centered / bordered style used by the existing
# index column (A2:A22) onto the newly appended rows.
$ws.Range("A23:A25").Font.Bold = $true
$ws.Range("A23:A25").HorizontalAlignment = -4108
$ws.Range("A23:A25").VerticalAlignment = -4160
$ws.Range("A23:A25").Borders.LineStyle = 1
